$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new "week" of data (2 rows: Primera/Segunda) right before the
# current top data block (old rows 373:374), pushing all subsequent rows
# down by 2 (old 375 -> new 377, ..., old 496 -> new 498).
#
# Copy the existing rows 373:374 and insert them above themselves; Excel's
# Insert-with-copied-clipboard pastes the duplicated content into the
# freshly inserted rows while shifting the original rows (and everything
# below) down by two. That gives us two new rows that already carry the
# correct static columns (A,B,C,E,F,G,H,I,N,O,Q,R), identical to the row
# that used to be there - we only need to overwrite the few cells that
# actually hold new data for the new week.
$ws.Rows("373:374").Copy()
$ws.Rows("373:374").Insert()

# New week's data for row 373 (Primera)
$ws.Range("D373").Value = 44985
$ws.Range("J373").Value = 200
$ws.Range("K373").Value = 3000
$ws.Range("L373").Value = 3500
$ws.Range("M373").Value = 3250
$ws.Range("P373").Value = 46

# New week's data for row 374 (Segunda)
$ws.Range("D374").Value = 44985
$ws.Range("J374").Value = 200
$ws.Range("K374").Value = 2500
$ws.Range("L374").Value = 3000
$ws.Range("M374").Value = 2750
$ws.Range("P374").Value = 28
